$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# "US in Tasks zerlegt": tasks 3-7 added for the 4. Userstory (Admin user
# management), plus "x" markers for already-finished tasks higher up the
# sheet, and two missing end-dates filled in.
#
# New shared strings must be introduced in this exact order so they line
# up with the ones already used elsewhere on the sheet:
#   3.Task / Es soll... / in Arbeit / x / 4.Task / Neben jeden User... /
#   Ein Admin soll... / 6. Task / Der User Klasse... / Wenn ein Benutzer...
# ---------------------------------------------------------------------

$ws.Range("C30").Value = "3.Task"
$ws.Range("D30").Value = "Es soll zwischen den verschiedenen Userrollen unterschieden werden "
$ws.Range("M32").Value = "in Arbeit "
$ws.Range("K20").Value = "x"
$ws.Range("C31").Value = "4.Task "
$ws.Range("D31").Value = "Neben jeden User in der Liste soll ein Button zum Löschen des Benutzers sein"
$ws.Range("D32").Value = "Ein Admin soll über ein Suchfeld nach der ID oder dem Nachnamen suchen können "
$ws.Range("C33").Value = "6. Task"
$ws.Range("D33").Value = "Der User Klasse soll ein Feld hinzugefügt werden ob ein Benutzer gesperrt ist"
$ws.Range("D34").Value = "Wenn ein Benutzer vom Admin als gesperrt erklärt wurde, soll er sich nicht mehr anmelden können "

# Mark the already-completed tasks higher up the sheet.
$ws.Range("K24").Value = "x"
$ws.Range("K25").Value = "x"

# End dates that were missing for the "1. Userstory" tasks (copy the date
# format off an existing date cell so we reuse its style instead of
# minting a new numFmt).
$ws.Range("L8").Copy()
$ws.Range("L28").PasteSpecial(-4122)
$ws.Range("L28").Value = 43516
$ws.Range("L29").PasteSpecial(-4122)
$ws.Range("L29").Value = 43516

# Remaining columns for the 5 new task rows (30-34) of "4. Userstory".
$ws.Range("C32").Value = "5. Task"
$ws.Range("C34").Value = "7. Task"

$ws.Range("K30").Value = "Manuel"
$ws.Range("K31").Value = "Simon"
$ws.Range("K32").Value = "Simon"
$ws.Range("K33").Value = "Simon"
$ws.Range("K34").Value = "Manuel"

$ws.Range("L30").PasteSpecial(-4122)
$ws.Range("L30").Value = 43522
$ws.Range("L31").PasteSpecial(-4122)
$ws.Range("L31").Value = 43522
$ws.Range("L32").PasteSpecial(-4122)
$ws.Range("L32").Value = 43522
$ws.Range("L33").PasteSpecial(-4122)
$ws.Range("L33").Value = 43522
$ws.Range("L34").PasteSpecial(-4122)
$ws.Range("L34").Value = 43522
$excel.CutCopyMode = $false

$ws.Range("M30").Value = "abgeschlossen"
$ws.Range("M31").Value = "abgeschlossen"
$ws.Range("M33").Value = "in Arbeit "

# Scroll the sheet back to the top (the saved view no longer pins
# topLeftCell to row 13) and reselect the same active cell as before.
$ws.Range("M30").Select()
$excel.ActiveWindow.ScrollRow = 1

Write-Output "done"
